$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select column N (DSSC_KEYWORDS) and delete the entire column, shifting
# everything to its right one column to the left - this removes the
# "keywords" field from the Performance model.
$ws.Range("N1:N1048576").Select()
$ws.Columns("N").Delete()
